$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.073692111274723
$ws.Range("D2").Value = 1.058218211496913
$ws.Range("E2").Value = 1.077075311459607
$ws.Range("F2").Value = 1.086562834020236
$ws.Range("I2").Value = 1.057206421313754
$ws.Range("J2").Value = 1.078605341527467
$ws.Range("K2").Value = 1.060951046734055
$ws.Range("L2").Value = 1.079757565694581
$ws.Range("M2").Value = 1.089220356331565
$ws.Range("N2").Value = 1.080137084830198
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.075292073299121
$ws.Range("D3").Value = 1.059070374326858
$ws.Range("E3").Value = 1.078523867392174
$ws.Range("F3").Value = 1.088099767924962
$ws.Range("I3").Value = 1.057712079487659
$ws.Range("J3").Value = 1.079861051147557
$ws.Range("K3").Value = 1.061617102541083
$ws.Range("L3").Value = 1.081022206009091
$ws.Range("M3").Value = 1.090574969929079
$ws.Range("N3").Value = 1.081394577701981
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.076325182139354
$ws.Range("D4").Value = 1.059620665597557
$ws.Range("E4").Value = 1.079458937069582
$ws.Range("F4").Value = 1.089092373374839
$ws.Range("I4").Value = 1.058037076411512
$ws.Range("J4").Value = 1.080670999157529
$ws.Range("K4").Value = 1.062046327911244
$ws.Range("L4").Value = 1.081837748434588
$ws.Range("M4").Value = 1.09144909363406
$ws.Range("N4").Value = 1.082205675931029
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.076758991419494
$ws.Range("D5").Value = 1.059851744365672
$ws.Range("E5").Value = 1.079851512805148
$ws.Range("F5").Value = 1.08950922094694
$ws.Range("I5").Value = 1.058173182423677
$ws.Range("J5").Value = 1.081010892097526
$ws.Range("K5").Value = 1.062226356927745
$ws.Range("L5").Value = 1.082179949408066
$ws.Range("M5").Value = 1.091816009101814
$ws.Range("N5").Value = 1.082546051557987
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.076831800327847
$ws.Range("D6").Value = 1.059890528139427
$ws.Range("E6").Value = 1.079917397334065
$ws.Range("F6").Value = 1.089579185742899
$ws.Range("I6").Value = 1.058196004670526
$ws.Range("J6").Value = 1.081067926138834
$ws.Range("K6").Value = 1.062256560216599
$ws.Range("L6").Value = 1.082237368414383
$ws.Range("M6").Value = 1.091877582881484
$ws.Range("N6").Value = 1.082603166594175
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.076330980709682
$ws.Range("D7").Value = 1.059623754315454
$ws.Range("E7").Value = 1.079464184747116
$ws.Range("F7").Value = 1.089097945044359
$ws.Range("I7").Value = 1.058038897116307
$ws.Range("J7").Value = 1.080675543206938
$ws.Range("K7").Value = 1.062048735103463
$ws.Range("L7").Value = 1.081842323496043
$ws.Range("M7").Value = 1.091453998590084
$ws.Range("N7").Value = 1.08221022643351
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.074233282436081
$ws.Range("D8").Value = 1.058506436772909
$ws.Range("E8").Value = 1.077565327792587
$ws.Range("F8").Value = 1.087082645878248
$ws.Range("I8").Value = 1.057377767997524
$ws.Range("J8").Value = 1.07903025381357
$ws.Range("K8").Value = 1.061176509303592
$ws.Range("L8").Value = 1.080185534483794
$ws.Range("M8").Value = 1.089678657023304
$ws.Range("N8").Value = 1.080562600540486
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.070519742466358
$ws.Range("D9").Value = 1.056528892533042
$ws.Range("E9").Value = 1.074201698882506
$ws.Range("F9").Value = 1.083516498515569
$ws.Range("I9").Value = 1.056195786632188
$ws.Range("J9").Value = 1.076110901907652
$ws.Range("K9").Value = 1.059625925244611
$ws.Range("L9").Value = 1.077244494794832
$ws.Range("M9").Value = 1.086531483393146
$ws.Range("N9").Value = 1.077639102820007
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.068031858937276
$ws.Range("D10").Value = 1.055204478347131
$ws.Range("E10").Value = 1.071946852621314
$ws.Range("F10").Value = 1.081128422552664
$ws.Range("I10").Value = 1.055396158548628
$ws.Range("J10").Value = 1.07415058916676
$ws.Range("K10").Value = 1.058582836191676
$ws.Range("L10").Value = 1.075268749280789
$ws.Range("M10").Value = 1.084420169773246
$ws.Range("N10").Value = 1.075676006210172
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.066951535703125
$ws.Range("D11").Value = 1.05462950937498
$ws.Range("E11").Value = 1.070967398969375
$ws.Range("F11").Value = 1.080091699306351
$ws.Range("I11").Value = 1.055047100472049
$ws.Range("J11").Value = 1.073298292953624
$ws.Range("K11").Value = 1.058128896889095
$ws.Range("L11").Value = 1.074409535722764
$ws.Range("M11").Value = 1.083502691673395
$ws.Range("N11").Value = 1.074822499638653
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.066549783757778
$ws.Range("D12").Value = 1.054415712399632
$ws.Range("E12").Value = 1.070603110454593
$ws.Range("F12").Value = 1.079706201163994
$ws.Range("I12").Value = 1.054917017691414
$ws.Range("J12").Value = 1.072981180754675
$ws.Range("K12").Value = 1.057959937569468
$ws.Range("L12").Value = 1.074089818512156
$ws.Range("M12").Value = 1.083161397957304
$ws.Range("N12").Value = 1.074504937104009
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.066635982474831
$ws.Range("D13").Value = 1.054461582957648
$ws.Range("E13").Value = 1.070681273295662
$ws.Range("F13").Value = 1.079788910701536
$ws.Range("I13").Value = 1.054944940285076
$ws.Range("J13").Value = 1.073049226614972
$ws.Range("K13").Value = 1.057996195645838
$ws.Range("L13").Value = 1.074158424772463
$ws.Range("M13").Value = 1.083234629488179
$ws.Range("N13").Value = 1.074573079597232
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.06691833646256
$ws.Range("D14").Value = 1.05461184152391
$ws.Range("E14").Value = 1.07093729654073
$ws.Range("F14").Value = 1.080059842402399
$ws.Range("I14").Value = 1.055036356526064
$ws.Range("J14").Value = 1.073272091262178
$ws.Range("K14").Value = 1.058114937743293
$ws.Range("L14").Value = 1.074383119428348
$ws.Range("M14").Value = 1.083474490522289
$ws.Range("N14").Value = 1.0747962607378
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.067092241162517
$ws.Range("D15").Value = 1.054704390460107
$ws.Range("E15").Value = 1.071094977514557
$ws.Range("F15").Value = 1.080226717289673
$ws.Range("I15").Value = 1.055092624368391
$ws.Range("J15").Value = 1.073409334792507
$ws.Range("K15").Value = 1.058188052684544
$ws.Range("L15").Value = 1.074521485760324
$ws.Range("M15").Value = 1.08362221004005
$ws.Range("N15").Value = 1.074933699169684
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.068103490990374
$ws.Range("D16").Value = 1.055242605416903
$ws.Range("E16").Value = 1.072011789642857
$ws.Range("F16").Value = 1.081197169194076
$ws.Range("I16").Value = 1.055419264731356
$ws.Range("J16").Value = 1.074207079307629
$ws.Range("K16").Value = 1.058612914355383
$ws.Range("L16").Value = 1.075325693527868
$ws.Range("M16").Value = 1.084480990046267
$ws.Range("N16").Value = 1.07573257657352
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.068736994653948
$ws.Range("D17").Value = 1.055579811866357
$ws.Range("E17").Value = 1.072586046190968
$ws.Range("F17").Value = 1.081805184785446
$ws.Range("I17").Value = 1.055623401219478
$ws.Range("J17").Value = 1.074706547588737
$ws.Range("K17").Value = 1.058878806806675
$ws.Range("L17").Value = 1.075829152935588
$ws.Range("M17").Value = 1.085018798150176
$ws.Range("N17").Value = 1.076232754156876
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.069106212918528
$ws.Range("D18").Value = 1.055776355314805
$ws.Range("E18").Value = 1.07292070293128
$ws.Range("F18").Value = 1.082159573118439
$ws.Range("I18").Value = 1.055742199379061
$ws.Range("J18").Value = 1.074997545098437
$ws.Range("K18").Value = 1.059033678176112
$ws.Range("L18").Value = 1.076122455763534
$ws.Range("M18").Value = 1.085332178200135
$ws.Range("N18").Value = 1.076524164916418
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.069232057381656
$ws.Range("D19").Value = 1.055843347360086
$ws.Range("E19").Value = 1.073034762136203
$ws.Range("F19").Value = 1.08228036716478
$ws.Range("I19").Value = 1.055782660611139
$ws.Range("J19").Value = 1.075096711409844
$ws.Range("K19").Value = 1.059086448259441
$ws.Range("L19").Value = 1.076222404298647
$ws.Range("M19").Value = 1.085438979713396
$ws.Range("N19").Value = 1.076623472055362
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.068669056180053
$ws.Range("D20").Value = 1.055543647665381
$ws.Range("E20").Value = 1.072524464741742
$ws.Range("F20").Value = 1.081739977121366
$ws.Range("I20").Value = 1.055601527389007
$ws.Range("J20").Value = 1.074652993980157
$ws.Range("K20").Value = 1.058850301766772
$ws.Range("L20").Value = 1.075775173460791
$ws.Range("M20").Value = 1.084961129020717
$ws.Range("N20").Value = 1.07617912449603
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.066835203380149
$ws.Range("D21").Value = 1.054567600446433
$ws.Range("E21").Value = 1.070861917304311
$ws.Range("F21").Value = 1.079980071264092
$ws.Range("I21").Value = 1.055009448545711
$ws.Range("J21").Value = 1.07320647793285
$ws.Range("K21").Value = 1.058079980742473
$ws.Range("L21").Value = 1.074316968192414
$ws.Range("M21").Value = 1.083403871358116
$ws.Range("N21").Value = 1.074730554230018
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.065679451113413
$ws.Range("D22").Value = 1.053952600145298
$ws.Range("E22").Value = 1.069813849017293
$ws.Range("F22").Value = 1.078871153395993
$ws.Range("I22").Value = 1.054634712013716
$ws.Range("J22").Value = 1.072293915913283
$ws.Range("K22").Value = 1.057593645040713
$ws.Range("L22").Value = 1.073396850556955
$ws.Range("M22").Value = 1.08242185465824
$ws.Range("N22").Value = 1.073816696267711
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.066292401017797
$ws.Range("D23").Value = 1.054278750044902
$ws.Range("E23").Value = 1.070369715218411
$ws.Range("F23").Value = 1.079459242725362
$ws.Range("I23").Value = 1.054833602811581
$ws.Range("J23").Value = 1.072777977759081
$ws.Range("K23").Value = 1.057851652223552
$ws.Range("L23").Value = 1.073884937423503
$ws.Range("M23").Value = 1.082942719511294
$ws.Range("N23").Value = 1.074301445536854
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.068699755559576
$ws.Range("D24").Value = 1.05555998915435
$ws.Range("E24").Value = 1.07255229166687
$ws.Range("F24").Value = 1.081769442451731
$ws.Range("I24").Value = 1.055611412069498
$ws.Range("J24").Value = 1.074677193582421
$ws.Range("K24").Value = 1.058863182644764
$ws.Range("L24").Value = 1.075799565563607
$ws.Range("M24").Value = 1.084987188189257
$ws.Range("N24").Value = 1.076203358464505
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.071481880360432
$ws.Range("D25").Value = 1.057041186147721
$ws.Range("E25").Value = 1.075073424097786
$ws.Range("F25").Value = 1.084440264598693
$ws.Range("I25").Value = 1.056503392664703
$ws.Range("J25").Value = 1.076868066136146
$ws.Range("K25").Value = 1.06002842334272
$ws.Range("L25").Value = 1.078007436318532
$ws.Range("M25").Value = 1.087347387983532
$ws.Range("N25").Value = 1.078397342308553
